# Add a new worksheet "sheet3" at the end of the workbook with example data,
# and make it the active (selected) sheet/tab.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it ends up last.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "sheet3"

# Populate the example data cells (see commit: "Add another example data sheet").
$ws3.Range("B2").Value = "[B2] start of ill-format excel file"
$ws3.Range("A8").Value = "[A8] second data appearance"

$ws3.Range("A13").Value = "[A13] contiguous column data"
$ws3.Range("B13").Value = "[B13] contiguous column data"
$ws3.Range("C13").Value = "[C13] contiguous column data"

$ws3.Range("A16").Value = "[A16] contiguous row data"
$ws3.Range("A17").Value = "[A17] contiguous row data"
$ws3.Range("A18").Value = "[A18] contiguous row data"

$ws3.Range("A22").Value = "[A2] jumping column data"
$ws3.Range("G22").Value = "[G22] jumping column data"
$ws3.Range("H22").Value = "[H22] jumping column data"

# Make the new sheet the active tab, scrolled to row 7 with J14 selected.
[void]$ws3.Activate()
[void]$ws3.Range("J14").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 7
$aw.ScrollColumn = 1
